$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.199.78'
$ws.Range('D3').Value = '1.859.09'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.70%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.32'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +3.62%  '
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.31'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +6.82%  '
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '2.128.46'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').Value = '1.844.40'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '35.158.30'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.92'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.19'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.30'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.88'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +24.44%  '
$ws.Range('E27').Value = '  +3.36%  '
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0561'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('B31').Value = 'BinanceUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.01'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('E33').Value = '  +27.95%  '
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('E35').Value = '  +10.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.812'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +16.93%  '
$ws.Range('E37').Value = '  +7.84%  '
$ws.Range('E38').Value = '  +4.29%  '
$ws.Range('E39').Value = '  +4.18%  '
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('D41').Value = '1.348.99'
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.05'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +3.26%  '
$ws.Range('E43').Value = '  +13.24%  '
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.24'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +42.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.56'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +4.76%  '
$ws.Range('D49').Value = '2.042.57'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0685'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('E51').Value = '  +0.68%  '
